# Fixes to generate proper choice list from csv and ajax queries
# (change from display.text to display.title.text)

$wb = $excel.ActiveWorkbook
$queries = $wb.Worksheets.Item("queries")

# Update the region-list callback (regions_csv row) to wrap the
# data_value/display output using display.title.text instead of display.text,
# and drop the now-unused name/label fields.
$regionCallback = "_.chain(context).pluck('region').uniq().map(function(region){`nreturn { data_value:region, display: {title: {text: region} } };`n}).value()"
$queries.Range("D2").Value = $regionCallback

# Update the country-list callback (countries_csv row) to return a much
# simpler object shape with data_value/region/display.title.text.
$countryCallback = "_.map(context, function(place){`nreturn { data_value: place.country, region: place.region, display: {title: {text: place.country} } };`n})"
$queries.Range("D3").Value = $countryCallback

# Make "queries" the active/selected sheet (was "settings"), with D8 selected.
[void]$queries.Activate()
[void]$queries.Range("D8").Select()
